# Updates cryptos list prices/percentages (and shifts the last three rows
# to insert a new coin), matching the "Updated cryptos list" GitHub Actions
# commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number (e.g. "216.26")
# must be forced back to text, otherwise Excel auto-converts the literal
# into a floating point number (and mangles the precision on save). We
# flip the cell to Text format, assign the literal, then restore the
# default "Normal" style so no stray number format lingers on the cell.
function Set-TextValue($rangeRef, $text) {
    $rng = $ws.Range($rangeRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# --- Price (column D) updates ---
$ws.Range("D2").Value = "26.726.57"
$ws.Range("D3").Value = "1.642.54"
Set-TextValue "D5" "216.26"
Set-TextValue "D10" "19.12"
Set-TextValue "D11" "0.0842"
$ws.Range("D12").Value = "1.866.77"
$ws.Range("D13").Value = "1.639.04"
Set-TextValue "D16" "64.33"
$ws.Range("D17").Value = "26.727.63"
$ws.Range("D18").Value = "0.0₃0734"
Set-TextValue "D19" "213.64"
Set-TextValue "D21" "4.37"
Set-TextValue "D22" "2.44"
Set-TextValue "D23" "6.24"
Set-TextValue "D24" "9.32"
Set-TextValue "D25" "145.37"
Set-TextValue "D28" "7.10"
Set-TextValue "D30" "0.0508"
Set-TextValue "D32" "3.35"
$ws.Range("D34").Value = "1.291.20"
Set-TextValue "D38" "0.533"
Set-TextValue "D41" "0.806"
Set-TextValue "D42" "2.24"
Set-TextValue "D43" "5.32"
$ws.Range("D44").Value = "1.792.65"
Set-TextValue "D45" "61.36"
Set-TextValue "D46" "91.22"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  -0.87%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("E19").Value = "  -1.88%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E22").Value = "  +12.71%  "
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("E24").Value = "  -2.40%  "
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("E27").Value = "  -1.49%  "
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("E37").Value = "  -3.15%  "
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("E41").Value = "  -0.63%  "
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("E43").Value = "  -2.78%  "
$ws.Range("E44").Value = "  +0.55%  "
$ws.Range("E45").Value = "  +2.95%  "
$ws.Range("E46").Value = "  -2.06%  "
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("E48").Value = "  +1.61%  "

# --- Row shift at the bottom of the list: a new coin (BabyDogeCoin) is
# inserted before EnergySwap, pushing EnergySwap down to replace Algorand,
# and Algorand down to replace USDD (which drops off the list). ---
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0102"
$ws.Range("E49").Value = "  -2.84%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D50" "7.69"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D51" "0.0973"
$ws.Range("E51").Value = "  -0.19%  "
